$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H88").Value = 1702.7059
$ws.Range("I88").Value = 1561.375
$ws.Range("K88").Value = 1561.375
$ws.Range("M88").Value = -1155.375
$ws.Range("H91").Value = 1702.7059
$ws.Range("I91").Value = 1561.375
$ws.Range("K91").Value = 1561.375
$ws.Range("M91").Value = -157.375
$ws.Range("H112").Value = 1061.7632
$ws.Range("J112").Value = 1074.2433
$ws.Range("L112").Value = 3222.7299
$ws.Range("N112").Value = -5438.7299
$ws.Range("H132").Value = 7872.636
$ws.Range("I132").Value = 9449.125
$ws.Range("J132").Value = 3668.6667
$ws.Range("K132").Value = 28347.375
$ws.Range("L132").Value = 11006.0001
$ws.Range("M132").Value = -25817.375
$ws.Range("N132").Value = -16066.0001
$ws.Range("H137").Value = 39381.594
$ws.Range("I137").Value = 2480.6667
$ws.Range("J137").Value = 113183.445
$ws.Range("K137").Value = 7442.000100000001
$ws.Range("L137").Value = 339550.335
$ws.Range("M137").Value = -4892.000100000001
$ws.Range("N137").Value = -344650.335
$ws.Range("H138").Value = 1378.41
$ws.Range("I138").Value = 549.93335
$ws.Range("J138").Value = 2056.2546
$ws.Range("K138").Value = 1649.80005
$ws.Range("L138").Value = 6168.763800000001
$ws.Range("M138").Value = 3490.19995
$ws.Range("N138").Value = -16448.7638

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1517.72
$ws.Range("I2").Value = 1306.55
$ws.Range("J2").Value = 2362.4
$ws.Range("K2").Value = 1306.55
$ws.Range("L2").Value = 2362.4
$ws.Range("M2").Value = -1193.55
$ws.Range("N2").Value = -2588.4
$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("N109").ClearContents()
$ws.Range("H116").Value = 1517.72
$ws.Range("I116").Value = 1306.55
$ws.Range("J116").Value = 2362.4
$ws.Range("K116").Value = 1306.55
$ws.Range("L116").Value = 2362.4
$ws.Range("M116").Value = 987.45
$ws.Range("N116").Value = -6950.4
$ws.Range("H132").Value = 15869.486
$ws.Range("I132").Value = 2051.7307
$ws.Range("K132").Value = 6155.1921
$ws.Range("M132").Value = -3625.1921

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1517.72
$ws.Range("I3").Value = 1306.55
$ws.Range("J3").Value = 2362.4
$ws.Range("K3").Value = 1306.55
$ws.Range("L3").Value = 2362.4
$ws.Range("M3").Value = -1192.55
$ws.Range("N3").Value = -2590.4
$ws.Range("H99").Value = 1332.2941
$ws.Range("I99").Value = 934.53845
$ws.Range("J99").Value = 2625
$ws.Range("K99").Value = 934.53845
$ws.Range("L99").Value = 2625
$ws.Range("M99").Value = 563.46155
$ws.Range("N99").Value = -5621
$ws.Range("H134").Value = 29120.514
$ws.Range("I134").Value = 37289.434
$ws.Range("K134").Value = 111868.302
$ws.Range("M134").Value = -109333.302

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 15191.909
$ws.Range("I31").Value = 28209.066
$ws.Range("K31").Value = 28209.066
$ws.Range("M31").Value = -27914.066
$ws.Range("H34").Value = 15191.909
$ws.Range("I34").Value = 28209.066
$ws.Range("K34").Value = 28209.066
$ws.Range("M34").Value = -28007.066
$ws.Range("H58").Value = 15375.343
$ws.Range("I58").Value = 1184.6086
$ws.Range("K58").Value = 1184.6086
$ws.Range("M58").Value = -981.6086
$ws.Range("H86").Value = 5754416
$ws.Range("I86").Value = 2512.4443
$ws.Range("K86").Value = 2512.4443
$ws.Range("M86").Value = -1389.4443
$ws.Range("H89").Value = 5754416
$ws.Range("I89").Value = 2512.4443
$ws.Range("K89").Value = 12562.2215
$ws.Range("M89").Value = -6946.2215
$ws.Range("H99").Value = 15155196
$ws.Range("I99").Value = 2994.15
$ws.Range("J99").Value = 38466276
$ws.Range("K99").Value = 2994.15
$ws.Range("L99").Value = 38466276
$ws.Range("M99").Value = -1496.15
$ws.Range("N99").Value = -38469272
$ws.Range("H122").Value = 1153.75
$ws.Range("I122").Value = 1106.7778
$ws.Range("J122").Value = 1214.1428
$ws.Range("K122").Value = 3320.3334
$ws.Range("L122").Value = 3642.4284
$ws.Range("M122").Value = -870.3334000000004
$ws.Range("N122").Value = -8542.428400000001
$ws.Range("H126").Value = 15155196
$ws.Range("I126").Value = 2994.15
$ws.Range("J126").Value = 38466276
$ws.Range("K126").Value = 8982.450000000001
$ws.Range("L126").Value = 115398828
$ws.Range("M126").Value = -6512.450000000001
$ws.Range("N126").Value = -115403768
$ws.Range("H136").Value = 15375.343
$ws.Range("I136").Value = 1184.6086
$ws.Range("K136").Value = 3553.8258
$ws.Range("M136").Value = -1003.8258

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 2351
$ws.Range("J39").Value = 2351
$ws.Range("L39").Value = 7053
$ws.Range("N39").Value = -7641
$ws.Range("H60").Value = 733.3333
$ws.Range("H113").Value = 613.2222
$ws.Range("I113").Value = 514.75
$ws.Range("J113").Value = 641.3570999999999
$ws.Range("K113").Value = 1544.25
$ws.Range("L113").Value = 1924.0713
$ws.Range("M113").Value = 625.75
$ws.Range("N113").Value = -6264.0713
$ws.Range("H131").Value = 775.5
$ws.Range("J131").Value = 779.8969
$ws.Range("L131").Value = 2339.6907
$ws.Range("N131").Value = -12419.6907
$ws.Range("H132").Value = 1474.3334
$ws.Range("I132").Value = 1435.8
$ws.Range("K132").Value = 12922.2
$ws.Range("M132").Value = -10392.2

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 83334430
$ws.Range("I122").Value = 30304094
$ws.Range("K122").Value = 90912282
$ws.Range("M122").Value = -90909832

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6800
$ws.Range("I7").Value = 6500
$ws.Range("J7").Value = 7700
$ws.Range("K7").Value = 6500
$ws.Range("L7").Value = 7700
$ws.Range("M7").Value = -6388
$ws.Range("N7").Value = -7924
$ws.Range("H22").Value = 1346.72
$ws.Range("I22").Value = 1086.75
$ws.Range("J22").Value = 1808.8889
$ws.Range("K22").Value = 1086.75
$ws.Range("L22").Value = 1808.8889
$ws.Range("M22").Value = -791.75
$ws.Range("N22").Value = -2398.8889
$ws.Range("H27").Value = 1346.72
$ws.Range("I27").Value = 1086.75
$ws.Range("J27").Value = 1808.8889
$ws.Range("K27").Value = 1086.75
$ws.Range("L27").Value = 1808.8889
$ws.Range("M27").Value = -979.75
$ws.Range("N27").Value = -2022.8889
$ws.Range("H40").Value = 3675.325
$ws.Range("I40").Value = 2281.5881
$ws.Range("J40").Value = 4705.478
$ws.Range("K40").Value = 2281.5881
$ws.Range("L40").Value = 4705.478
$ws.Range("M40").Value = -2145.5881
$ws.Range("N40").Value = -4977.478
$ws.Range("H46").Value = 999.7857
$ws.Range("I46").Value = 762.125
$ws.Range("J46").Value = 1316.6666
$ws.Range("K46").Value = 762.125
$ws.Range("L46").Value = 1316.6666
$ws.Range("M46").Value = -574.125
$ws.Range("N46").Value = -1692.6666
$ws.Range("H126").Value = 6800
$ws.Range("I126").Value = 6500
$ws.Range("J126").Value = 7700
$ws.Range("K126").Value = 19500
$ws.Range("L126").Value = 23100
$ws.Range("M126").Value = -17030
$ws.Range("N126").Value = -28040
$ws.Range("H127").Value = 0
$ws.Range("J127").Value = 0
$ws.Range("L127").Value = 0
$ws.Range("N127").ClearContents()
$ws.Range("H128").Value = 0
$ws.Range("J128").Value = 0
$ws.Range("L128").Value = 0
$ws.Range("N128").ClearContents()

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H69").Value = 10890.4
$ws.Range("J69").Value = 12612.75
$ws.Range("L69").Value = 12612.75
$ws.Range("N69").Value = -14110.75
$ws.Range("H72").Value = 10890.4
$ws.Range("J72").Value = 12612.75
$ws.Range("L72").Value = 37838.25
$ws.Range("N72").Value = -45326.25
$ws.Range("H122").Value = 1696.4242
$ws.Range("I122").Value = 1697.4348
$ws.Range("J122").Value = 1694.1
$ws.Range("K122").Value = 5092.3044
$ws.Range("L122").Value = 5082.299999999999
$ws.Range("M122").Value = -2642.3044
$ws.Range("N122").Value = -9982.299999999999
$ws.Range("H132").Value = 1141.2
$ws.Range("I132").Value = 897.6818
$ws.Range("K132").Value = 2693.0454
$ws.Range("M132").Value = -163.0454
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()
